$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Append text to the review comment (D17, task 14 observation)
$old = $ws.Range("D17").Value2
$ws.Range("D17").Value2 = $old + ".  A mayores se añadieron comentarios para entender lo que se hacía"

# Clear the "Estado" cell for task 14 (row 17) - status could not be set to "corregido"
$ws.Range("C17").Value2 = ""

# Set "Tiempo empleado" for task 14 (row 17, column G)
$ws.Range("G17").Value2 = 33

# Reset D17 cell style to the plain "Estado" style (remove the special centered/wrap font)
$ws.Range("C17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Adjust row height to fit new (longer) text
$ws.Rows.Item(17).RowHeight = 256.7
